# Apply roster reorder/update to the active worksheet (rows 2-15 change,
# row 12 and rows 16-19 stay the same, per the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jordan Clarkson", "SG,SF", "Utah Jazz"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("AJ Green", "PG,SG", "Milwaukee Bucks"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Alperen Sengun", "C", "Houston Rockets"),
    @("Nic Claxton", "C", "Brooklyn Nets"),
    @("Jared McCain", "PG", "Philadelphia 76ers"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Goga Bitadze", "C", "Orlando Magic"),
    @("Devin Booker", "PG,SG", "Phoenix Suns")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
